# Swap the two worker records that sit in rows 16 and 17 of Hoja1.
# Previously:
#   Row16: CC 73190466 GERMAN ENRIQUE DUPERRET TRESPALACIOS 2001  33125    828116
#   Row17: CC 45757793 YENIS ESTHER CASTELLAR CASTELLAR      2201  40000   1000000
# After the edit the two data rows trade places (values only, formatting stays
# anchored to the row it was already on):
#   Row16: CC 45757793 YENIS ESTHER CASTELLAR CASTELLAR      2201  40000   1000000
#   Row17: CC 73190466 GERMAN ENRIQUE DUPERRET TRESPALACIOS 2001  33125    828116

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$cols = @("C", "D", "E", "F", "G")

$row16Values = @{}
$row17Values = @{}

foreach ($col in $cols) {
    $row16Values[$col] = $ws.Range($col + "16").Value2
    $row17Values[$col] = $ws.Range($col + "17").Value2
}

foreach ($col in $cols) {
    $ws.Range($col + "16").Value2 = $row17Values[$col]
    $ws.Range($col + "17").Value2 = $row16Values[$col]
}
